$wb = $excel.ActiveWorkbook

# --- Step 1: rename the existing (only) sheet to "no_adjust" ---
$base = $wb.Worksheets.Item(1)
$base.Name = "no_adjust"

# --- Step 2: duplicate it, placing the copy BEFORE the original, and rename to "adjust" ---
$base.Copy($base) | Out-Null
$adjust = $wb.Worksheets.Item(1)
$adjust.Name = "adjust"

# --- Step 3: on the "adjust" sheet, turn row 4 (B:P) into formulas referencing row 5 ---
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
foreach ($col in $cols) {
    $adjust.Range("$col`4").Formula = "=0.306*" + $col + "5"
}

# Apply the new number format / alignment style (numFmtId 2 "0.00", horizontal left)
# to the B4:P13 block (column A keeps the default style).
# Build the combined format on a scratch cell first and paste only the
# formatting across - setting NumberFormat then HorizontalAlignment
# directly on the destination range would materialize two separate
# intermediate styles instead of a single combined one.
$helper = $adjust.Range("Z1")
$helper.NumberFormat = "0.00"
$helper.HorizontalAlignment = -4131  # xlLeft
$helper.Copy()
$styleRange = $adjust.Range("B4:P13")
$styleRange.PasteSpecial(-4122)  # xlPasteFormats
$helper.Clear()

# --- Step 4: sheet view / selection tweaks ---
$noAdjust = $wb.Worksheets.Item("no_adjust")
$noAdjust.Activate()
$noAdjust.Range("C15").Select()

$adjust.Activate()
$adjust.Range("B4").Select()
